$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Artn"
$ws.Range("C2").Value = "Gfra1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.089056333333333
$ws.Range("H2").Value = 6.267169
$ws.Range("I2").Value = 0.7196603919224289
$ws.Range("J2").Value = 0.719660391922429
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.06547366666666667
$ws.Range("N2").Value = 0.196421
$ws.Range("O2").Value = 0.002125877360986814
$ws.Range("P2").Value = 0.002125877360986814
$ws.Range("Q2").Value = 0.1367781780165555
$ws.Range("R2").Value = 1.231003602149
$ws.Range("S2").Value = 0.00152990973478679
$ws.Range("T2").Value = 0.00152990973478679

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Artn"
$ws.Range("C3").Value = "Gfra1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.089056333333333
$ws.Range("H3").Value = 6.267169
$ws.Range("I3").Value = 0.7196603919224289
$ws.Range("J3").Value = 0.719660391922429
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 25.94643066666667
$ws.Range("N3").Value = 77.839292
$ws.Range("O3").Value = 0.8424597607080814
$ws.Range("P3").Value = 0.8424597607080814
$ws.Range("Q3").Value = 54.20355531159422
$ws.Range("R3").Value = 487.831997804348
$ws.Range("S3").Value = 0.6062849215700535
$ws.Range("T3").Value = 0.6062849215700536

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Artn"
$ws.Range("C4").Value = "Gfra1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.089056333333333
$ws.Range("H4").Value = 6.267169
$ws.Range("I4").Value = 0.7196603919224289
$ws.Range("J4").Value = 0.719660391922429
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.786517
$ws.Range("N4").Value = 14.359551
$ws.Range("O4").Value = 0.1554143619309319
$ws.Range("P4").Value = 0.1554143619309319
$ws.Range("Q4").Value = 9.999303653457666
$ws.Range("R4").Value = 89.993732881119
$ws.Range("S4").Value = 0.1118455606175886
$ws.Range("T4").Value = 0.1118455606175887

$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Artn"
$ws.Range("C5").Value = "Gfra1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.8137799999999999
$ws.Range("H5").Value = 2.44134
$ws.Range("I5").Value = 0.280339608077571
$ws.Range("J5").Value = 0.280339608077571
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.06547366666666667
$ws.Range("N5").Value = 0.196421
$ws.Range("O5").Value = 0.002125877360986814
$ws.Range("P5").Value = 0.002125877360986814
$ws.Range("Q5").Value = 0.05328116046
$ws.Range("R5").Value = 0.47953044414
$ws.Range("S5").Value = 0.0005959676262000244
$ws.Range("T5").Value = 0.0005959676262000244

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Artn"
$ws.Range("C6").Value = "Gfra1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.8137799999999999
$ws.Range("H6").Value = 2.44134
$ws.Range("I6").Value = 0.280339608077571
$ws.Range("J6").Value = 0.280339608077571
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 25.94643066666667
$ws.Range("N6").Value = 77.839292
$ws.Range("O6").Value = 0.8424597607080814
$ws.Range("P6").Value = 0.8424597607080814
$ws.Range("Q6").Value = 21.11468634792
$ws.Range("R6").Value = 190.03217713128
$ws.Range("S6").Value = 0.2361748391380278
$ws.Range("T6").Value = 0.2361748391380278

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Artn"
$ws.Range("C7").Value = "Gfra1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.8137799999999999
$ws.Range("H7").Value = 2.44134
$ws.Range("I7").Value = 0.280339608077571
$ws.Range("J7").Value = 0.280339608077571
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.786517
$ws.Range("N7").Value = 14.359551
$ws.Range("O7").Value = 0.1554143619309319
$ws.Range("P7").Value = 0.1554143619309319
$ws.Range("Q7").Value = 3.89517180426
$ws.Range("R7").Value = 35.05654623834
$ws.Range("S7").Value = 0.04356880131334321
$ws.Range("T7").Value = 0.04356880131334321
